$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "kalimat" (sentence) column (column B) contains Indonesian example
# questions that refer to "barang" (goods). Rename all occurrences of
# "barang" to "produk" to match the renamed dataset column/table, per the
# commit message ("Fix metric inside function evaluate"). Column C (the
# "label"/answer column, e.g. "nama_barang, stok") must stay untouched.
$lastRow = $ws.UsedRange.Rows.Count
$colB = $ws.Range("B1:B$lastRow")
$colB.Replace("barang", "produk", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart) | Out-Null

# Restore the last-selected cell as reflected by the saved workbook.
$ws.Range("C3").Select() | Out-Null
